# Daily cryptos-list refresh (Coinranking snapshot), matching the GitHub Actions commit
# "Updated cryptos list on Fri Jul 28 04:26:25 UTC 2023 with GitHub Actions".
#
# Column D ("Price") cells are plain-text numeric-looking strings (e.g. "1.001"), exactly
# as the source sheet stores them (t="inlineStr"). Typing a bare numeric-looking literal into
# Range.Value lets Excel auto-coerce it to a real number, so for any new Price value that
# would otherwise parse as a number we prefix it with a leading apostrophe (Excel's own
# "force text" convention) to keep it text, matching the source data. Values that already
# contain two dots (e.g. "29.226.24") can never parse as a number, so no prefix is needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.226.24"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
$ws.Range("D3").Value = "1.861.69"
$ws.Range("E3").Value = "  -0.90%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "'0.7134"
$ws.Range("E5").Value = "  -0.90%  "

# Row 6
$ws.Range("D6").Value = "'240.63"
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 -> Cardano
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3082"
$ws.Range("E8").Value = "  -1.12%  "

# Row 9 -> Dogecoin
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.07714"
$ws.Range("E9").Value = "  -1.37%  "

# Row 10
$ws.Range("D10").Value = "'24.87"
$ws.Range("E10").Value = "  -0.79%  "

# Row 11
$ws.Range("D11").Value = "'0.08316"
$ws.Range("E11").Value = "  +0.81%  "

# Row 12
$ws.Range("D12").Value = "1.878.57"
$ws.Range("E12").Value = "  -0.65%  "

# Row 13
$ws.Range("D13").Value = "'0.7143"
$ws.Range("E13").Value = "  -1.99%  "

# Row 14
$ws.Range("D14").Value = "'5.201"
$ws.Range("E14").Value = "  -1.79%  "

# Row 15
$ws.Range("D15").Value = "'90.84"
$ws.Range("E15").Value = "  -0.67%  "

# Row 16
$ws.Range("D16").Value = "29.266.22"
$ws.Range("E16").Value = "  -1.07%  "

# Row 17
$ws.Range("D17").Value = "'5.956"
$ws.Range("E17").Value = "  +0.07%  "

# Row 18
$ws.Range("D18").Value = "'242.56"
$ws.Range("E18").Value = "  -1.64%  "

# Row 19
$ws.Range("D19").Value = "'0.000007810"
$ws.Range("E19").Value = "  -1.05%  "

# Row 20
$ws.Range("D20").Value = "2.130.56"
$ws.Range("E20").Value = "  +0.11%  "

# Row 21
$ws.Range("D21").Value = "'13.15"
$ws.Range("E21").Value = "  -1.22%  "

# Row 22
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").Value = "'7.894"
$ws.Range("E23").Value = "  -0.92%  "

# Row 24
$ws.Range("D24").Value = "'0.9998"
$ws.Range("E24").Value = "  -0.04%  "

# Row 25
$ws.Range("E25").Value = "  +1.35%  "

# Row 26
$ws.Range("D26").Value = "'163.24"
$ws.Range("E26").Value = "  -0.39%  "

# Row 27
$ws.Range("D27").Value = "'8.884"
$ws.Range("E27").Value = "  -1.82%  "

# Row 28
$ws.Range("D28").Value = "'18.55"
$ws.Range("E28").Value = "  +1.12%  "

# Row 29
$ws.Range("D29").Value = "'1.349"
$ws.Range("E29").Value = "  -1.17%  "

# Row 30
$ws.Range("D30").Value = "'1.500"
$ws.Range("E30").Value = "  +1.13%  "

# Row 31
$ws.Range("D31").Value = "'4.418"
$ws.Range("E31").Value = "  +0.65%  "

# Row 32
$ws.Range("D32").Value = "'4.254"
$ws.Range("E32").Value = "  +2.50%  "

# Row 33
$ws.Range("D33").Value = "'0.05159"
$ws.Range("E33").Value = "  -2.33%  "

# Row 34
$ws.Range("D34").Value = "'0.8300"
$ws.Range("E34").Value = "  +14.87%  "

# Row 35
$ws.Range("D35").Value = "'1.930"
$ws.Range("E35").Value = "  -0.94%  "

# Row 36
$ws.Range("D36").Value = "'1.170"
$ws.Range("E36").Value = "  -2.72%  "

# Row 37
$ws.Range("D37").Value = "'2.683"
$ws.Range("E37").Value = "  +0.25%  "

# Row 38
$ws.Range("E38").Value = "  -0.47%  "

# Row 39
$ws.Range("D39").Value = "'2.691"
$ws.Range("E39").Value = "  -1.24%  "

# Row 40
$ws.Range("D40").Value = "1.161.56"
$ws.Range("E40").Value = "  -5.84%  "

# Row 41
$ws.Range("D41").Value = "'6.193"
$ws.Range("E41").Value = "  +1.32%  "

# Row 42
$ws.Range("D42").Value = "'0.8957"
$ws.Range("E42").Value = "  -1.27%  "

# Row 43
$ws.Range("D43").Value = "'72.77"
$ws.Range("E43").Value = "  -1.49%  "

# Row 44
$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("D45").Value = "'101.94"
$ws.Range("E45").Value = "  -1.51%  "

# Row 46
$ws.Range("D46").Value = "2.029.27"
$ws.Range("E46").Value = "  -0.14%  "

# Row 47
$ws.Range("D47").Value = "'0.5168"
$ws.Range("E47").Value = "  -3.14%  "

# Row 48
$ws.Range("D48").Value = "'1.785"
$ws.Range("E48").Value = "  +1.05%  "

# Row 49
$ws.Range("D49").Value = "'9.333"
$ws.Range("E49").Value = "  +0.47%  "

# Row 50 -> Frax
$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  +0.04%  "

# Row 51 -> Aptos
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'7.051"
$ws.Range("E51").Value = "  -0.52%  "
